# The deck's date placeholders ("datetimeFigureOut" fields) cache their
# last-rendered value as literal text inside the field. PowerPoint
# refreshes this cache whenever the file is saved. Here the cached date
# moves from 09.04.2025 to 10.04.2025 on the Slide Master and on every
# Slide Layout (the placeholder is inherited, but each layout keeps its
# own cached copy of the field text).

$newDate = "10.04.2025"
$ppPlaceholderDate = 16

$p = $ppt.ActivePresentation
$m = $p.SlideMaster

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePlaceholder = $false
        if ($sh.HasTextFrame) {
            try {
                if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                    $isDatePlaceholder = $true
                }
            } catch {
                $isDatePlaceholder = $false
            }
        }
        if ($isDatePlaceholder) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide Master date placeholder.
Update-DatePlaceholder $m.Shapes

# Every Slide Layout keeps its own cached date placeholder text too.
for ($li = 1; $li -le $m.CustomLayouts.Count; $li++) {
    $layout = $m.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}
